# koersen bel 20 s3
# Append another snapshot block (rows 106:126, the most recent poll) to the
# bottom of the sheet as rows 127:147 - mirrors the "fetch every 10s / write
# CSV / ship to S3" job re-emitting the same reading.
#
# Copy+PasteSpecial(values) is used instead of `.Value = .Value` because the
# source cells are text that merely look numeric (e.g. "79,850"); a plain
# Value round-trip gets reinterpreted as a number and also stamps a new
# NumberFormat style on the destination cells, neither of which happened in
# the real edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A106:F126")
$target = $ws.Range("A127:F147")

$source.Copy()
$target.PasteSpecial(-4163)  # xlPasteValues - keep cell text/type, no formatting carried over

$excel.CutCopyMode = $false
